# Update the "Correspond Handoff Datetime" (D2) and
# "Correspond Handback DateTime" (G2) timestamps on the zh-cn and de-de
# report sheets to reflect a newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-18 03:48:43"
$wsZhCn.Range("G2").Value = "2016-01-18 03:49:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-18 03:48:54"
$wsDeDe.Range("G2").Value = "2016-01-18 03:49:46"
